# Add season-record columns (Wins / Losses / Ties) to the team stats sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of an existing header cell (bold, centered, bordered)
# onto the three new header cells so they match the rest of row 1.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the season record for every data row (2-54) with the team's
# win/loss/tie totals.
$lastRow = 54
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 95
    $ws.Cells.Item($r, 31).Value = 67
    $ws.Cells.Item($r, 32).Value = 0
}
